$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 gets new data (previously blank) ---
$ws.Range("A4").Value = "secretsauce.com/:track/:page/"
$ws.Range("B4").Value = "expectations"

# --- Alignment updates ---
# Column A "label" cells (A1:A3, A5) center aligned
foreach ($addr in @("A1","A2","A3","A5")) {
  $ws.Range($addr).HorizontalAlignment = -4108
}

# Column B "value" cells (B1:B4) left aligned
foreach ($addr in @("B1","B2","B3","B4")) {
  $ws.Range($addr).HorizontalAlignment = -4131
}

# A4 (new url row) right aligned
$ws.Range("A4").HorizontalAlignment = -4152

# Header row 5 (B5 and the A/B/C/D letter headers C5:F5) center aligned
foreach ($addr in @("B5","C5","D5","E5","F5")) {
  $ws.Range($addr).HorizontalAlignment = -4108
}

# --- Update selected cell to B4 ---
$ws.Range("B4").Select() | Out-Null
